# chore: update Sheets via scheduled runner
# Refreshes cached market-price-derived figures (currentAveragePrice*,
# LevePrice*, LeveProfit*) across the per-job Leve-profit sheets.
$wb = $excel.ActiveWorkbook

# ALC row 17
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 638039
$ws.Range("J17").Value = 638039
$ws.Range("L17").Value = 1914117
$ws.Range("N17").Value = -1914453

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2620.3845
$ws.Range("I132").Value = 2620.3845
$ws.Range("K132").Value = 7861.1535
$ws.Range("M132").Value = -5331.1535

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2604.818
$ws.Range("I137").Value = 2165.3
$ws.Range("K137").Value = 6495.900000000001
$ws.Range("M137").Value = -3945.900000000001

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1990.8823
$ws.Range("I138").Value = 1729.8667
$ws.Range("J138").Value = 3948.5
$ws.Range("K138").Value = 5189.6001
$ws.Range("L138").Value = 11845.5
$ws.Range("M138").Value = -49.60009999999966
$ws.Range("N138").Value = -22125.5

# ARM row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 701509.25
$ws.Range("I2").Value = 981445.7
$ws.Range("K2").Value = 981445.7
$ws.Range("M2").Value = -981332.7

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 125001870
$ws.Range("I61").Value = 125001870
$ws.Range("K61").Value = 125001870
$ws.Range("M61").Value = -125001658

# ARM row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 701509.25
$ws.Range("I116").Value = 981445.7
$ws.Range("K116").Value = 981445.7
$ws.Range("M116").Value = -979151.7

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2440702.8
$ws.Range("I132").Value = 2501595.2
$ws.Range("K132").Value = 7504785.600000001
$ws.Range("M132").Value = -7502255.600000001

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 125001870
$ws.Range("I136").Value = 125001870
$ws.Range("K136").Value = 375005610
$ws.Range("M136").Value = -375003060

# BSM row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 701509.25
$ws.Range("I3").Value = 981445.7
$ws.Range("K3").Value = 981445.7
$ws.Range("M3").Value = -981331.7

# BSM row 62
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()

# BSM row 65
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()

# BSM row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2171.9688
$ws.Range("I86").Value = 2028.6
$ws.Range("K86").Value = 2028.6
$ws.Range("M86").Value = -905.5999999999999

# BSM row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 2171.9688
$ws.Range("I89").Value = 2028.6
$ws.Range("K89").Value = 10143
$ws.Range("M89").Value = -4527

# BSM row 92
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H92").Value = 30000
$ws.Range("J92").Value = 30000
$ws.Range("L92").Value = 30000
$ws.Range("N92").Value = -34992

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 100001120
$ws.Range("I134").Value = 125001050
$ws.Range("K134").Value = 375003150
$ws.Range("M134").Value = -375000615

# CRP row 7
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 305.125
$ws.Range("I7").Value = 382.5
$ws.Range("J7").Value = 73
$ws.Range("K7").Value = 382.5
$ws.Range("L7").Value = 73
$ws.Range("M7").Value = -269.5
$ws.Range("N7").Value = -299

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10619.533
$ws.Range("I31").Value = 4286.625
$ws.Range("J31").Value = 17857.143
$ws.Range("K31").Value = 4286.625
$ws.Range("L31").Value = 17857.143
$ws.Range("M31").Value = -3991.625
$ws.Range("N31").Value = -18447.143

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 10619.533
$ws.Range("I34").Value = 4286.625
$ws.Range("J34").Value = 17857.143
$ws.Range("K34").Value = 4286.625
$ws.Range("L34").Value = 17857.143
$ws.Range("M34").Value = -4084.625
$ws.Range("N34").Value = -18261.143

# CRP row 62
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2750
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 2750
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 2750
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -3998

# CRP row 65
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 2750
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 2750
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 13750
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -19990

# CUL row 113
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 72695.21000000001
$ws.Range("J113").Value = 1516.75
$ws.Range("L113").Value = 4550.25
$ws.Range("N113").Value = -8890.25

# CUL row 122
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1157.2
$ws.Range("I122").Value = 973
$ws.Range("J122").Value = 1894
$ws.Range("K122").Value = 8757
$ws.Range("L122").Value = 17046
$ws.Range("M122").Value = -6307
$ws.Range("N122").Value = -21946

# GSM row 68
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H68").Value = 54999
$ws.Range("J68").Value = 54999
$ws.Range("L68").Value = 54999
$ws.Range("N68").Value = -56621

# GSM row 71
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H71").Value = 54999
$ws.Range("J71").Value = 54999
$ws.Range("L71").Value = 164997
$ws.Range("N71").Value = -173109

# GSM row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2508.6072
$ws.Range("I102").Value = 2571.88
$ws.Range("J102").Value = 1981.3334
$ws.Range("K102").Value = 2571.88
$ws.Range("L102").Value = 1981.3334
$ws.Range("M102").Value = -949.8800000000001
$ws.Range("N102").Value = -5225.3334

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 5684438.5
$ws.Range("I132").Value = 6946829.5
$ws.Range("J132").Value = 3679.5
$ws.Range("K132").Value = 20840488.5
$ws.Range("L132").Value = 11038.5
$ws.Range("M132").Value = -20837958.5
$ws.Range("N132").Value = -16098.5

# LTW row 22
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3593.125
$ws.Range("J22").Value = 2498.3333
$ws.Range("L22").Value = 2498.3333
$ws.Range("N22").Value = -3088.3333

# LTW row 27
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 3593.125
$ws.Range("J27").Value = 2498.3333
$ws.Range("L27").Value = 2498.3333
$ws.Range("N27").Value = -2712.3333

# LTW row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2607.1428
$ws.Range("I46").Value = 2625
$ws.Range("K46").Value = 2625
$ws.Range("M46").Value = -2437

# LTW row 68
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 5279157
$ws.Range("J68").Value = 19999
$ws.Range("L68").Value = 19999
$ws.Range("N68").Value = -21497

# LTW row 71
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 5279157
$ws.Range("J71").Value = 19999
$ws.Range("L71").Value = 99995
$ws.Range("N71").Value = -107483

# WVR row 126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2226.318
$ws.Range("I126").Value = 2321.923
$ws.Range("J126").Value = 2088.2222
$ws.Range("K126").Value = 6965.768999999999
$ws.Range("L126").Value = 6264.6666
$ws.Range("M126").Value = -4495.768999999999
$ws.Range("N126").Value = -11204.6666

# WVR row 130
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H130").Value = 72499.75
$ws.Range("J130").Value = 72499.75
$ws.Range("L130").Value = 72499.75
$ws.Range("N130").Value = -82539.75

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 10002715
$ws.Range("J132").Value = 7362
$ws.Range("L132").Value = 22086
$ws.Range("N132").Value = -27146

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 29413296
$ws.Range("I136").Value = 29413296
$ws.Range("K136").Value = 88239888
$ws.Range("M136").Value = -88237338
